$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '36.980.27'
$ws.Range("E2").Value = '  -1.11%  '

# Row 3
$ws.Range("D3").Value = '2.008.25'
$ws.Range("E3").Value = '  -1.98%  '

# Row 4
$ws.Range("E4").Value = '  -0.09%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '226.00'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.11%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.603'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.05%  '

# Row 7
$ws.Range("E7").Value = '  -0.03%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '54.98'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.15%  '

# Row 9
$ws.Range("E9").Value = '  -3.46%  '

# Row 10
$ws.Range("E10").Value = '  -3.86%  '

# Row 11
$ws.Range("E11").Value = '  -4.11%  '

# Row 12
$ws.Range("D12").Value = '2.306.04'
$ws.Range("E12").Value = '  -1.87%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '13.99'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.54%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '19.67'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.36%  '

# Row 15
$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.14'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.82%  '

# Row 16
$ws.Range("B16").Value = 'Polygon'
$ws.Range("C16").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.731'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.77%  '

# Row 17
$ws.Range("D17").Value = '2.008.17'
$ws.Range("E17").Value = '  -2.09%  '

# Row 18
$ws.Range("D18").Value = '36.905.32'
$ws.Range("E18").Value = '  -1.08%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.18'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.64%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '68.54'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.68%  '

# Row 21
$ws.Range("D21").Value = '0.0₃0809'
$ws.Range("E21").Value = '  -3.97%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '221.33'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.92%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.06%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.44'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.26%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.17'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.81%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '164.17'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.58%  '

# Row 27
$ws.Range("E27").Value = '  -6.68%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.127'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.13%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.53'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.92%  '

# Row 30
$ws.Range("E30").Value = '  -5.57%  '

# Row 31
$ws.Range("E31").Value = '  -1.17%  '

# Row 32
$ws.Range("E32").Value = '  -3.12%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0598'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.16%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.42'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.50%  '

# Row 35
$ws.Range("E35").Value = '  -3.07%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.87'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.33%  '

# Row 37
$ws.Range("E37").Value = '  -0.09%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.13'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.56%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.34'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.63%  '

# Row 40
$ws.Range("D40").Value = '1.461.96'
$ws.Range("E40").Value = '  -2.67%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0210'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.18%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '94.12'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.86%  '

# Row 43
$ws.Range("B43").Value = 'FTXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.23'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +10.14%  '

# Row 44
$ws.Range("B44").Value = 'Cronos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0909'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.49%  '

# Row 45
$ws.Range("B45").Value = 'HuobiToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.76'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.16%  '

# Row 46
$ws.Range("E46").Value = '  -2.02%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '15.73'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -6.36%  '

# Row 48
$ws.Range("E48").Value = '  -2.14%  '

# Row 49
$ws.Range("E49").Value = '  -1.93%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.89'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.13%  '

# Row 51
$ws.Range("D51").Value = '2.193.15'
$ws.Range("E51").Value = '  -1.89%  '
